$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.950.93"
$ws.Range("E2").Value = "  -2.15%  "

$ws.Range("D3").Value = "1.982.62"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.017"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.015"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4897"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4146"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08791"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.080"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.13%  "

$ws.Range("D12").Value = "2.069.91"
$ws.Range("E12").Value = "  +2.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.867"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.351"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.06%  "

$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001095"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06662"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.015"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.924"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("D23").Value = "29.005.01"
$ws.Range("E23").Value = "  -2.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("D26").Value = "2.286.01"
$ws.Range("E26").Value = "  +1.33%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.74%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.167"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.210"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.028"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09792"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.502"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.791"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "

$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02386"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.299"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06323"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.934"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6414"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.90%  "

$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.014"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1956"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.352"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6128"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000344"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.486"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.167"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.32%  "
